$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Conversion en tick" test block (G2:M8) ---
# L4: alim threshold test value
$ws.Range("L4").Value = 93
# H5: wheel diameter test value
$ws.Range("H5").Value = 193.55

# --- New "Coeff en" block (rows 12-14, columns A-B) ---
$ws.Range("A12").Value = "Tension alim avec test"
$ws.Range("B12").Value = 13
$ws.Range("A13").Value = "Tension en cours"
$ws.Range("B13").Value = 10
$ws.Range("A14").Value = "Coeff en"
$ws.Range("B14").Formula = "=B12/B13"

# --- Update the precision table (K15:K17) ---
$ws.Range("K15").Value = 100.28
$ws.Range("K16").Value = 101.269
$ws.Range("K17").Value = 90.09

# --- Update sheet view / selection ---
$ws.Range("G16").Select() | Out-Null
